# Update the "民族自治地方国内生产总值" sheet so it only keeps the
# 2010年-2013年 rows (previously rows 8-11), moved up to rows 2-5,
# and the old rows for 2000年/2005年-2009年 are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the data currently sitting in rows 8-11 (years 2010-2013)
# before we start overwriting anything.
$sourceRows = 8..11
$data = @{}
foreach ($r in $sourceRows) {
    $rowVals = @{}
    foreach ($c in 1..6) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $data[$r] = $rowVals
}

# Write that captured data into rows 2-5.
$destRow = 2
foreach ($r in $sourceRows) {
    foreach ($c in 1..6) {
        $ws.Cells.Item($destRow, $c).Value2 = $data[$r][$c]
    }
    $destRow = $destRow + 1
}

# Delete the now-obsolete rows 6-11 (old 2008/2009 duplicates plus the
# original 2010-2013 rows that were copied up above).
$ws.Range("A6:F11").EntireRow.Delete() | Out-Null
